$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.11"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "1.788.08"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'307.94"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4536"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").Value = "'0.3690"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'0.07223"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "'0.8519"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").Value = "'20.34"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "1.783.30"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.07012"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.459"
$ws.Range("E15").Value = "  -4.16%  "
$ws.Range("D16").Value = "'90.28"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'0.000008583"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "26.661.27"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'5.251"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").Value = "2.017.04"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'1.903"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").Value = "'149.58"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.05"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.129"
$ws.Range("E28").Value = "  -12.42%  "
$ws.Range("D29").Value = "'5.181"
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("D30").Value = "'113.44"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").Value = "'0.08823"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'0.7502"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'4.428"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("D35").Value = "'2.876"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'1.108"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "'0.01938"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "'0.05199"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "'2.874"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "'7.104"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.323"
$ws.Range("E42").Value = "  +5.82%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5173"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'0.1637"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").Value = "'8.444"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "'0.4922"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "'103.51"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'1.638"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "'0.06270"
$ws.Range("E51").Value = "  -1.60%  "

# Values that look like plain numbers (e.g. "0.3690", "1.001") were entered
# with a leading apostrophe above so Excel keeps them as literal text
# (matching the source data, which stores them as inline strings) instead of
# silently re-parsing them into doubles and dropping trailing zeros. That
# leading apostrophe flips the cell's style to "quote prefixed", so copy the
# plain/default formatting from a known plain-text cell (B2) back over the
# whole edited block to restore the original (unstyled) appearance without
# touching any of the values we just wrote.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D2:E51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
